# Fix the bug that was sending two different passwords for the same
# account: row 2 ("arrennbaral@gmail.com") had an unrelated password
# ("djgjhsgksdgjo") instead of the actual password used elsewhere for
# that login ("2kWip@HgY!S9gAL", the same rich-text value/hyperlink
# already present on B6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b2 = $ws.Range("B2")
$b6 = $ws.Range("B6")

# Give B2 a hyperlink (mailto:2kWip@HgY) matching the one already on B6 -
# Hyperlinks.Add() also happens to stamp the cell with the built-in
# "Hyperlink" style/text, which we'll immediately overwrite below.
$ws.Hyperlinks.Add($b2, "mailto:2kWip@HgY", "", "", "2kWip@HgY")

# Now copy B6's actual rich-text content + formatting (style, fonts,
# colors for the "2kWip@HgY" + "!S9gAL" runs) onto B2, so B2 ends up
# with the exact same password value/formatting as B6, instead of the
# stray plain-text hyperlink label.
$b6.Copy($b2)

# Move the active selection to B8, matching the saved workbook view.
$ws.Range("B8").Select()
